$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3921.658785588634
$ws.Range("B3").Value = 4610.862094733957
$ws.Range("B4").Value = 5293.444360066154
$ws.Range("B5").Value = 5992.858100300337
$ws.Range("B6").Value = 6648.827874564511
$ws.Range("B7").Value = 7163.696237356986
$ws.Range("B8").Value = 7405.689504784748
$ws.Range("B9").Value = 8165.13799976693
$ws.Range("B10").Value = 8538.713075064825
$ws.Range("B11").Value = 9058.50130673449
$ws.Range("B12").Value = 9336.660358264509
$ws.Range("B13").Value = 9599.174680571881
$ws.Range("B14").Value = 10007.49005345957
$ws.Range("B15").Value = 10321.62533764913
$ws.Range("B16").Value = 10596.27984208827
$ws.Range("B17").Value = 10814.84726068009
$ws.Range("B18").Value = 11138.68398215823
$ws.Range("B19").Value = 11386.5604879413
$ws.Range("B20").Value = 11296.02019712908
$ws.Range("B21").Value = 11482.99371222917
$ws.Range("B22").Value = 11676.5549147708
$ws.Range("B23").Value = 11826.54117721885
$ws.Range("B24").Value = 12158.21408126541
$ws.Range("B25").Value = 12370.79548344
$ws.Range("B26").Value = 12568.38345922562
$ws.Range("B27").Value = 12701.27060576151
$ws.Range("B28").Value = 12764.90919545843
$ws.Range("B29").Value = 12907.95048228432
$ws.Range("B30").Value = 13084.51181512031
$ws.Range("B31").Value = 13243.12588675456
$ws.Range("B32").Value = 13576.59158557177
$ws.Range("B33").Value = 13669.3144584906
$ws.Range("B34").Value = 13878.21174268121
$ws.Range("B35").Value = 14047.12327388834
$ws.Range("B36").Value = 14152.66769723202
$ws.Range("B37").Value = 14271.97539074888
$ws.Range("B38").Value = 14434.42331817191
$ws.Range("B39").Value = 14479.19064378691
$ws.Range("B40").Value = 14668.98115302313
$ws.Range("B41").Value = 14764.03071671831
$ws.Range("B42").Value = 14885.69482890222
$ws.Range("B43").Value = 14980.07477207137
$ws.Range("B44").Value = 14934.35330678715
$ws.Range("B45").Value = 15071.79478457906
$ws.Range("B46").Value = 15132.62758229976
$ws.Range("B47").Value = 15233.73239740377
$ws.Range("B48").Value = 15372.4660136289
$ws.Range("B49").Value = 15369.68032229482
$ws.Range("B50").Value = 15407.38825627885
$ws.Range("B51").Value = 15536.51856861041
$ws.Range("B52").Value = 15558.20613624086
$ws.Range("B53").Value = 15675.38133453084
$ws.Range("B54").Value = 15721.38662937669
$ws.Range("B55").Value = 15893.93215097395
$ws.Range("B56").Value = 15831.68663770271
$ws.Range("B57").Value = 15951.09853415263
$ws.Range("B58").Value = 15960.4713648951
$ws.Range("B59").Value = 15905.0524474193
$ws.Range("B60").Value = 16036.57893001658
$ws.Range("B61").Value = 16050.71011181273
$ws.Range("B62").Value = 16220.6029797299
